# b5-t3-p12 writing practice:
# Fix the typo "edtion" -> "edition" in the second paragraph by inserting
# the missing "i", producing three separate runs:
#   "This is the second ed" | "i" | "tion."

$d = $word.ActiveDocument

# Locate the typo text robustly (position-independent).
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("edtion", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'edtion' in the document."
}

$matchStart = $findRange.Start

# Insertion point is right after "ed", before "tion".
$insertPos = $matchStart + 2
$insertionRange = $d.Range($insertPos, $insertPos)
$insertionRange.InsertAfter("i")

# The newly inserted "i" character as its own range.
$newCharRange = $d.Range($insertPos, $insertPos + 1)

# Adding then immediately deleting a bookmark around the inserted
# character forces Word to keep it as a distinct run (split from its
# neighbors) without leaving any run-formatting residue behind.
$d.Bookmarks.Add("tmp_split_mark", $newCharRange)
$d.Bookmarks("tmp_split_mark").Delete()
